# Insert a new historical data row at row 483 (2019-11-21) into Sheet1,
# pushing the existing rows 483..555 down to 484..556.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 483.
$ws.Rows.Item(483).Insert()

# Populate the newly inserted row with the historical data point.
$ws.Range("A483").Value = 1574294400

# B and C look like a date / a zero-padded number respectively, so force
# them to be stored as plain text (matching the rest of the column),
# then clear the temporary Text number format so no stray cell style is
# left behind.
$ws.Range("B483").NumberFormat = "@"
$ws.Range("B483").Value = "2019-11-21"
$ws.Range("B483").ClearFormats()

$ws.Range("C483").NumberFormat = "@"
$ws.Range("C483").Value = "03001"
$ws.Range("C483").ClearFormats()

$ws.Range("D483").Value = "CLOUD"
$ws.Range("E483").Value = 0.165
$ws.Range("F483").Value = 0.165
$ws.Range("G483").Value = 0.165
$ws.Range("H483").Value = 0.165
$ws.Range("I483").Value = 50000
